# "Fruta / hortaliza, semanal" weekly refresh of the price series.
# Rows 5-21 get new Fecha / Volumen / Precio values (the weekly roll-forward),
# and two additional observations are appended as rows 22 and 23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant values shared by every data row in this sheet (needed for the two
# brand-new rows appended at the bottom).
$mercadoId = 10
$mercado   = 'Vega Modelo de Temuco'
$region    = 'La Araucanía'
$codreg    = 9
$categoriaId = 100112036
$categoria = 'Caigua'
$variedad  = 'Sin especificar'
$calidad   = 'Primera'
$unidad    = '$/caja 15 kilos'
$origen    = 'Región de Arica y Parinacota'
$kgUnidades = 15
$clasificacion = 'Hortaliza'
$dateFormat = 'YYYY-MM-DD HH:MM:SS'

# New D / J / K / L / M / P values, keyed by destination row.
$rows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23)
$fecha   = @{5=44838;6=44518;7=44757;8=44812;9=44825;10=44819;11=44813;12=44525;13=44756;14=44830;15=44771;16=44837;17=44776;18=44749;19=44811;20=44824;21=44755;22=44826;23=44767}
$volumen = @{5=10;6=50;7=30;8=80;9=30;10=100;11=20;12=40;13=80;14=25;15=40;16=80;17=80;18=50;19=30;20=20;21=50;22=50;23=50}
$pmin    = @{5=20000;6=10000;7=20000;8=20000;9=20000;10=20000;11=20000;12=8000;13=20000;14=12000;15=20000;16=16000;17=20000;18=20000;19=20000;20=20000;21=20000;22=20000;23=20000}
$pmax    = @{5=20000;6=10000;7=20000;8=20000;9=20000;10=20000;11=20000;12=8000;13=20000;14=12000;15=20000;16=16000;17=20000;18=20000;19=20000;20=20000;21=20000;22=20000;23=20000}
$pprom   = @{5=20000;6=10000;7=20000;8=20000;9=20000;10=20000;11=20000;12=8000;13=20000;14=12000;15=20000;16=16000;17=20000;18=20000;19=20000;20=20000;21=20000;22=20000;23=20000}
$pkg     = @{5=1333;6=667;7=1333;8=1333;9=1333;10=1333;11=1333;12=533;13=1333;14=800;15=1333;16=1067;17=1333;18=1333;19=1333;20=1333;21=1333;22=1333;23=1333}

foreach ($r in $rows) {
    if ($r -gt 21) {
        # Rows 22 & 23 are brand new - populate every column.
        $ws.Cells.Item($r, 1).Value = $mercadoId
        $ws.Cells.Item($r, 2).Value = $mercado
        $ws.Cells.Item($r, 3).Value = $region
        $ws.Cells.Item($r, 5).Value = $codreg
        $ws.Cells.Item($r, 6).Value = $categoriaId
        $ws.Cells.Item($r, 7).Value = $categoria
        $ws.Cells.Item($r, 8).Value = $variedad
        $ws.Cells.Item($r, 9).Value = $calidad
        $ws.Cells.Item($r, 14).Value = $unidad
        $ws.Cells.Item($r, 15).Value = $origen
        $ws.Cells.Item($r, 17).Value = $kgUnidades
        $ws.Cells.Item($r, 18).Value = $clasificacion
    }

    $dCell = $ws.Cells.Item($r, 4)
    $dCell.Value = $fecha[$r]
    $dCell.NumberFormat = $dateFormat

    $ws.Cells.Item($r, 10).Value = $volumen[$r]
    $ws.Cells.Item($r, 11).Value = $pmin[$r]
    $ws.Cells.Item($r, 12).Value = $pmax[$r]
    $ws.Cells.Item($r, 13).Value = $pprom[$r]
    $ws.Cells.Item($r, 16).Value = $pkg[$r]
}
